$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The R5 potentiometer spec (row 15) is being updated from the 20k-pot part
# to the 10k-pot part (the values that used to live in the now-redundant
# "R5b" row 16). REFDES stays "R5".
$ws.Range("B15").Value = "10k pot"
$ws.Range("D15").Value = "Bourns, inc"
$ws.Range("E15").Value = "PTV09A-4225F-B103"
$ws.Range("F15").Value = "PTV09A-4225F-B103-ND"

# Row 16 ("R5b") is now a redundant duplicate of the updated row 15, so its
# contents are cleared out entirely.
$ws.Range("A16:I16").ClearContents()

# H7's price was missing the shared currency number format that the rest of
# the Price column uses; bring it in line with its neighbours.
$ws.Range("H7").NumberFormat = $ws.Range("H8").NumberFormat

Write-Output "edit applied"
